# Update "liste des requetes.xlsx" (Feuil1) to reflect the removal of the
# REALMS_SEL_IPTEMPORAIRE_RECHERCHEIP query row and refresh of the
# Table-concernee AutoFilter, per issue #131 (maj de la base / ip_temporaire).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Make sure we're on the right sheet / clear any pre-existing AutoFilter so
# it gets rebuilt cleanly against the post-delete row count.
$ws.Select()
$ws.AutoFilterMode = $false

# Remove the now-obsolete "Savoir si une ip a deja eu une erreur
# d'authentification" / REALMS_SEL_IPTEMPORAIRE_RECHERCHEIP row (old row 28).
# Everything below shifts up by one (old 29 -> 28, old 30 -> 29).
$ws.Rows.Item(28).Delete()

# Re-apply the "Table concernee" AutoFilter over the whole table, now
# restricted to the ip_banned / ip_temporaire rows (this both redraws the
# dropdown arrows and hides every row whose column A isn't one of the two
# values).
$rng = $ws.UsedRange
$rng.AutoFilter(1, @("ip_banned", "ip_temporaire"), 7)

# Keep the workbook-level _FilterDatabase defined name in sync with the new
# data extent (was Feuil1!$A$1:$E$24, now covers the full 29 rows).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Feuil1!_FilterDatabase") {
        $n.RefersTo = "=Feuil1!`$A`$1:`$E`$29"
    }
}

# Restore the selection to where the author left off (C28, the Nom requete
# cell of the last visible row) instead of the old scrolled/selected state.
$ws.Range("C28").Select()
